# Auto-generated Excel COM-interop script applying the Raiden_Profits diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1965.25
$ws.Range("I2").Value = 1462.1666
$ws.Range("K2").Value = 1462.1666
$ws.Range("M2").Value = -1349.1666
$ws.Range("H15").Value = 785.6531
$ws.Range("I15").Value = 785.6531
$ws.Range("K15").Value = 2356.9593
$ws.Range("M15").Value = -2187.9593
$ws.Range("H17").Value = 1038.0952
$ws.Range("J17").Value = 1038.0952
$ws.Range("L17").Value = 3114.2856
$ws.Range("N17").Value = -3450.2856
$ws.Range("H19").Value = 1053.7222
$ws.Range("I19").Value = 994.4167
$ws.Range("J19").Value = 1172.3334
$ws.Range("K19").Value = 994.4167
$ws.Range("L19").Value = 1172.3334
$ws.Range("M19").Value = -819.4167
$ws.Range("N19").Value = -1522.3334
$ws.Range("H28").Value = 3691.4
$ws.Range("I28").Value = 3267
$ws.Range("K28").Value = 3267
$ws.Range("M28").Value = -2782
$ws.Range("H38").Value = 1945
$ws.Range("J38").Value = 3200
$ws.Range("L38").Value = 9600
$ws.Range("N38").Value = -10344
$ws.Range("H41").Value = 58
$ws.Range("I41").Value = 58
$ws.Range("K41").Value = 58
$ws.Range("M41").Value = 382
$ws.Range("H76").Value = 5166
$ws.Range("I76").Value = 5749
$ws.Range("K76").Value = 5749
$ws.Range("M76").Value = -5434
$ws.Range("H79").Value = 5166
$ws.Range("I79").Value = 5749
$ws.Range("K79").Value = 5749
$ws.Range("M79").Value = -4657
$ws.Range("H125").Value = 2180.1052
$ws.Range("I125").Value = 4487.6665
$ws.Range("J125").Value = 1115.0769
$ws.Range("K125").Value = 40388.9985
$ws.Range("L125").Value = 10035.6921
$ws.Range("M125").Value = -37928.9985
$ws.Range("N125").Value = -14955.6921
$ws.Range("H129").Value = 1949.7858
$ws.Range("I129").Value = 2574.25
$ws.Range("K129").Value = 7722.75
$ws.Range("M129").Value = -2722.75
$ws.Range("H134").Value = 67499.5
$ws.Range("J134").Value = 67499.5
$ws.Range("L134").Value = 67499.5
$ws.Range("N134").Value = -77639.5
$ws.Range("H138").Value = 4741.1523
$ws.Range("J138").Value = 5463.4683
$ws.Range("L138").Value = 16390.4049
$ws.Range("N138").Value = -26670.4049

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3999.5
$ws.Range("I2").Value = 3999
$ws.Range("K2").Value = 3999
$ws.Range("M2").Value = -3886
$ws.Range("H32").Value = 3152.0435
$ws.Range("I32").Value = 2084.7424
$ws.Range("K32").Value = 2084.7424
$ws.Range("M32").Value = -1797.7424
$ws.Range("H61").Value = 3497.4443
$ws.Range("I61").Value = 3497.4443
$ws.Range("K61").Value = 3497.4443
$ws.Range("M61").Value = -3285.4443
$ws.Range("H116").Value = 3999.5
$ws.Range("I116").Value = 3999
$ws.Range("K116").Value = 3999
$ws.Range("M116").Value = -1705
$ws.Range("H122").Value = 4034.3462
$ws.Range("I122").Value = 4025.3333
$ws.Range("K122").Value = 12075.9999
$ws.Range("M122").Value = -9625.999899999999
$ws.Range("H132").Value = 1858.1818
$ws.Range("I132").Value = 1844.05
$ws.Range("K132").Value = 5532.15
$ws.Range("M132").Value = -3002.15
$ws.Range("H136").Value = 3497.4443
$ws.Range("I136").Value = 3497.4443
$ws.Range("K136").Value = 10492.3329
$ws.Range("M136").Value = -7942.332900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3999.5
$ws.Range("I3").Value = 3999
$ws.Range("K3").Value = 3999
$ws.Range("M3").Value = -3885
$ws.Range("H9").Value = 450000
$ws.Range("J9").Value = 450000
$ws.Range("L9").Value = 450000
$ws.Range("N9").Value = -450336
$ws.Range("H64").Value = 997.25
$ws.Range("I64").Value = 997.25
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 997.25
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -772.25
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 997.25
$ws.Range("I67").Value = 997.25
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 997.25
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -217.25
$ws.Range("N67").ClearContents()
$ws.Range("H134").Value = 1881.6316
$ws.Range("I134").Value = 1732.2142
$ws.Range("K134").Value = 5196.642599999999
$ws.Range("M134").Value = -2661.642599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 434.22726
$ws.Range("I22").Value = 451.1579
$ws.Range("J22").Value = 327
$ws.Range("K22").Value = 451.1579
$ws.Range("L22").Value = 327
$ws.Range("M22").Value = -101.1579
$ws.Range("N22").Value = -1027
$ws.Range("H70").Value = 57925
$ws.Range("J70").Value = 57925
$ws.Range("L70").Value = 57925
$ws.Range("N70").Value = -58555
$ws.Range("H73").Value = 57925
$ws.Range("J73").Value = 57925
$ws.Range("L73").Value = 57925
$ws.Range("N73").Value = -60109
$ws.Range("H94").Value = 2335.6667
$ws.Range("I94").Value = 2253.5
$ws.Range("K94").Value = 2253.5
$ws.Range("M94").Value = -1802.5
$ws.Range("H132").Value = 1180.8182
$ws.Range("I132").Value = 1180.8182
$ws.Range("K132").Value = 3542.4546
$ws.Range("M132").Value = -1012.4546
$ws.Range("H134").Value = 1931.5385
$ws.Range("I134").Value = 1954.6111
$ws.Range("K134").Value = 5863.8333
$ws.Range("M134").Value = -3328.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 73399.2
$ws.Range("J37").Value = 73399.2
$ws.Range("L37").Value = 220197.6
$ws.Range("N37").Value = -220421.6
$ws.Range("H108").Value = 3357.7144
$ws.Range("I108").Value = 3357.7144
$ws.Range("K108").Value = 10073.1432
$ws.Range("M108").Value = -7193.143199999999
$ws.Range("H132").Value = 1169.1
$ws.Range("J132").Value = 1316
$ws.Range("L132").Value = 11844
$ws.Range("N132").Value = -16904
$ws.Range("H138").Value = 3264.75
$ws.Range("I138").Value = 1019.6667
$ws.Range("K138").Value = 3059.0001
$ws.Range("M138").Value = 2080.9999
$ws.Range("H139").Value = 8958.333000000001
$ws.Range("I139").Value = 3493
$ws.Range("J139").Value = 19889
$ws.Range("K139").Value = 10479
$ws.Range("L139").Value = 59667
$ws.Range("M139").Value = -5339
$ws.Range("N139").Value = -69947

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8539.608
$ws.Range("I70").Value = 7764.4375
$ws.Range("J70").Value = 10311.429
$ws.Range("K70").Value = 7764.4375
$ws.Range("L70").Value = 10311.429
$ws.Range("M70").Value = -7494.4375
$ws.Range("N70").Value = -10851.429
$ws.Range("H73").Value = 8539.608
$ws.Range("I73").Value = 7764.4375
$ws.Range("J73").Value = 10311.429
$ws.Range("K73").Value = 7764.4375
$ws.Range("L73").Value = 10311.429
$ws.Range("M73").Value = -6828.4375
$ws.Range("N73").Value = -12183.429
$ws.Range("H122").Value = 2676.5386
$ws.Range("J122").Value = 2999.3333
$ws.Range("L122").Value = 8997.999899999999
$ws.Range("N122").Value = -13897.9999
$ws.Range("H132").Value = 1873.9375
$ws.Range("I132").Value = 1732.8667
$ws.Range("K132").Value = 5198.6001
$ws.Range("M132").Value = -2668.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 791.75
$ws.Range("I16").Value = 500.13043
$ws.Range("K16").Value = 500.13043
$ws.Range("M16").Value = -330.13043
$ws.Range("H38").Value = 42999.668
$ws.Range("I38").Value = 49999
$ws.Range("K38").Value = 49999
$ws.Range("M38").Value = -49589
$ws.Range("H40").Value = 3557.8333
$ws.Range("I40").Value = 3529.2083
$ws.Range("K40").Value = 3529.2083
$ws.Range("M40").Value = -3393.2083
$ws.Range("H93").Value = 3913
$ws.Range("I93").Value = 5556
$ws.Range("K93").Value = 5556
$ws.Range("M93").Value = -4308
$ws.Range("H100").Value = 2642.7778
$ws.Range("I100").Value = 2196.25
$ws.Range("K100").Value = 2196.25
$ws.Range("M100").Value = -1655.25
$ws.Range("H122").Value = 3264.2856
$ws.Range("I122").Value = 3264.2856
$ws.Range("K122").Value = 9792.856800000001
$ws.Range("M122").Value = -7342.856800000001
$ws.Range("H132").Value = 2197.476
$ws.Range("I132").Value = 2134.0527
$ws.Range("K132").Value = 6402.158100000001
$ws.Range("M132").Value = -3872.158100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9524.75
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 9524.75
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H132").Value = 3041.2
$ws.Range("I132").Value = 3153.8096
$ws.Range("K132").Value = 9461.4288
$ws.Range("M132").Value = -6931.4288
$ws.Range("H136").Value = 706.61536
$ws.Range("J136").Value = 1222.8572
$ws.Range("L136").Value = 3668.5716
$ws.Range("N136").Value = -8768.571599999999

